{"js": "// Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// \"\u00a9 2020 ...\" copyright paragraph that used to follow the requirements\n// list (right after \"LOM3215: F\u00edsica do Estado S\u00f3lido (Requisito)\"),\n// while leaving the trailing empty paragraph (and the page-break\n// paragraph after it) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOM3215: ...\") so we only touch the\n// three paragraphs that immediately follow it.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOM3215\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the 'LOM3215' anchor paragraph.\");\n}\n\nconst targetTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst toDelete = [];\nfor (let offset = 0; offset < targetTexts.length; offset++) {\n  const idx = anchorIndex + 1 + offset;\n  if (idx >= items.length) {\n    break;\n  }\n  if (items[idx].text !== targetTexts[offset]) {\n    throw new Error(\n      \"Unexpected paragraph at index \" + idx + \": \" + JSON.stringify(items[idx].text)\n    );\n  }\n  toDelete.push(items[idx]);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n# \"\u00a9 2020 ...\" copyright paragraph that used to follow the requirements\n# list (right after \"LOM3215: F\u00edsica do Estado S\u00f3lido (Requisito)\"),\n# while leaving the trailing empty paragraph (and the page-break\n# paragraph after it) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"LOM3215: ...\") so we only touch the\n# three paragraphs that immediately follow it.\n$anchor = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOM3215*\") {\n        $anchor = $i\n        break\n    }\n}\n\nif ($anchor -eq -1) {\n    throw \"Could not locate the 'LOM3215' anchor paragraph.\"\n}\n\n# Expected text (paragraph mark stripped) of the three paragraphs that\n# must be removed, in document order right after the anchor.\n$targetTexts = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Verify the paragraphs are what we expect before deleting anything.\nfor ($k = 0; $k -lt $targetTexts.Length; $k++) {\n    $idx = $anchor + 1 + $k\n    $p = $d.Paragraphs.Item($idx)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -ne $targetTexts[$k]) {\n        throw \"Unexpected paragraph at index $idx : [$text]\"\n    }\n}\n\n# Delete the fixed paragraph right after the anchor, three times, so the\n# collection keeps re-settling and indices stay valid.\nfor ($k = 0; $k -lt $targetTexts.Length; $k++) {\n    $idx = $anchor + 1\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
